$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("John", "Joey", "Greg", "Arielie", "Lisa", "Bob")

for ($i = 0; $i -lt 6; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    for ($j = 0; $j -lt $names.Length; $j++) {
        $ws.Cells.Item($row, 2 + $j).Value = $names[$j]
    }
}
